$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row after the current last data row (row 33), so the
# existing row 33 content can be pushed down to row 34 while row 33
# receives the new, more recent price entry.
$ws.Rows.Item(34).Insert()

# Copy the (old) row 33 values down into the newly inserted row 34.
for ($col = 1; $col -le 20; $col++) {
    $src = $ws.Cells.Item(33, $col)
    $dst = $ws.Cells.Item(34, $col)
    $dst.Value2 = $src.Value2
}

# Overwrite row 33 with the new weekly price record.
$ws.Range("D33").Value2 = 44516
$ws.Range("E33").Value2 = 10
$ws.Range("F33").Value2 = "Fruta"
$ws.Range("G33").Value2 = 100103
$ws.Range("H33").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I33").Value2 = 100103001
$ws.Range("J33").Value2 = "Cereza"
$ws.Range("K33").Value2 = "Early Burlat"
$ws.Range("L33").Value2 = "Segunda"
$ws.Range("M33").Value2 = 600
$ws.Range("N33").Value2 = 13000
$ws.Range("O33").Value2 = 13500
$ws.Range("P33").Value2 = 13250
$ws.Range("Q33").Value2 = "`$/bandeja 5 kilos"
$ws.Range("R33").Value2 = "Provincia de Curicó"
$ws.Range("S33").Value2 = 2650
$ws.Range("T33").Value2 = 5
